$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add Wins / Losses / Ties in AD1:AF1.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the formatting of the existing header cell (AC1 - bold, bordered,
# centered/top-aligned) onto the new header cells so the new columns match
# the look of the rest of the header row, reusing the same style record.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Data rows 2-48: team record values (constant across all players on the roster).
$lastRow = 48
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 74  # AD: Wins
    $ws.Cells.Item($r, 31).Value = 88  # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF: Ties
}
